$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 in I1 and IF in J1, matching the H1 header style
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data rows 2-24: column I is always 1, column J mirrors column H
for ($r = 2; $r -le 24; $r++) {
    $hValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hValue
}
